# The deck's "last modified/printed" date placeholder (master + every
# slide layout) was bumped from 29-6-2021 to 25-10-2021, and the title
# slide's subtitle (which spells the date out as free text) was bumped
# from "1 July 2021" to "2 November 2021".

$p = $ppt.ActivePresentation

function Update-DatePlaceholders($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.Name -like "Date Placeholder*") {
            $shape.TextFrame.TextRange.Text = $newText
        }
    }
}

$newDate = "25-10-2021"

# Slide master's own date placeholder.
$master = $p.SlideMaster
Update-DatePlaceholders $master.Shapes $newDate

# Every slide layout has its own date placeholder too.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholders $layout.Shapes $newDate
}

# Title slide spells the date out in its subtitle text.
$titleSlide = $p.Slides.Item(1)
$titleSlide.Shapes.Item(2).TextFrame.TextRange.Text = "2 November 2021"
